$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the text in B14 into two parts: keep first sentence in B14, move the rest to C14
$ws.Range("B14").Value = "4. Scroll down the information to find the traceroutes. "
$ws.Range("C14").Value = " Traceroutes should be located near the bottom of the scan showing you a list of Ip addresses and the name of routers that the packet information is passing through."

# Update row 14 height
$ws.Rows("14").RowHeight = 44.4

# Update selection to B13
$ws.Range("B13").Select() | Out-Null
